# ------------------------------------------------------------------
# "spreadsheet as valuesets, added testcases, refactor"
#
# Starting layout (2 cols: CODE / DESCRIPTIVE_TEXT header + 3 sample
# rows A2..A4), ending layout (3 cols, with a new "options" section on
# top of the table plus a couple of extra example rows at the bottom).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert a new column C ---------------------------
$ws.Columns.Item(3).Insert()

# --- 2. Insert the 2 new "options" rows at the top; old row1 (CODE /
#        DESCRIPTIVE_TEXT header) becomes row3, etc. This also grows
#        the sheet from 23 to 25 rows (the 2 extra rows needed at the
#        bottom for the new trailing blank rows). -------------------
$ws.Rows.Item(1).Resize(2).Insert()

# ------------------------------------------------------------------
# Values
# ------------------------------------------------------------------
$ws.Range("A1").Value = "options.lookupOrder"
$ws.Range("B1").Value = "CODE"
$ws.Range("C1").Value = "DESCRIPTIVE_TEXT"

$ws.Range("A2").Value = "options.separator"
$ws.Range("B2").Value = "DESCRIPTIVE_TEXT"
$ws.Range("C2").Value = ";"

$ws.Range("A3").Value = "CODE"
$ws.Range("B3").Value = "DESCRIPTIVE_TEXT"
$ws.Range("C3").Value = "'"

$ws.Range("A4").Value = "A2:CODE"
$ws.Range("B4").Value = "B2:DESCRIPTIVE_TEXT"
$ws.Range("C4").Value = "'"

$ws.Range("A5").Value = "A3:CODE"
$ws.Range("B5").Value = "B3:DESCRIPTIVE_TEXT"
$ws.Range("C5").Value = "'"

$ws.Range("A6").Value = "A4"
$ws.Range("B6").Value = "X; Y; Z"
$ws.Range("C6").Value = "'"

# ------------------------------------------------------------------
# Column widths (closest reachable values - COM ColumnWidth is
# quantized to 1/6 character increments in this engine, so the exact
# fractional widths from the source file can only be approximated)
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 28.1
$ws.Columns.Item(2).ColumnWidth = 24.6
$ws.Columns.Item(3).ColumnWidth = 29.1

# Default (no-data) look for columns A & C vs column B, matching the
# base style each column falls back to outside the populated rows.
$ws.Columns.Item(1).HorizontalAlignment = -4131
$ws.Columns.Item(1).WrapText = $true
$ws.Columns.Item(3).HorizontalAlignment = -4131
$ws.Columns.Item(3).WrapText = $true

# ------------------------------------------------------------------
# Row heights
# ------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 17.25
$ws.Rows.Item(2).RowHeight = 17.25
$ws.Rows.Item(3).RowHeight = 27.649999999999995
$ws.Rows.Item(4).RowHeight = 20.45
for ($r = 5; $r -le 18; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}
for ($r = 19; $r -le 25; $r++) {
    $ws.Rows.Item($r).RowHeight = 20.25
}

# ------------------------------------------------------------------
# Formatting - options rows (1 & 2): plain default font, wrap on A/C,
# general alignment + no-wrap on B
# ------------------------------------------------------------------
$optsRows = $ws.Range("A1:C2")
$optsRows.WrapText = $true
$optsRows.HorizontalAlignment = -4131
$ws.Range("B1:B2").WrapText = $false
$ws.Range("B1:B2").HorizontalAlignment = 1

# ------------------------------------------------------------------
# Formatting - header + sample rows (3..6): black font, wrap text
# ------------------------------------------------------------------
$dataRows = $ws.Range("A3:C6")
$dataRows.Font.Color = 0
$dataRows.WrapText = $true
$dataRows.HorizontalAlignment = -4131

# header row 3 all left aligned (already set above)
# sample rows 4 & 5: column B & C right-aligned
$ws.Range("B4:B5").HorizontalAlignment = -4152
$ws.Range("C4:C5").HorizontalAlignment = -4152

# ------------------------------------------------------------------
# Formatting - trailing blank rows (7..25)
# ------------------------------------------------------------------
$blankA = $ws.Range("A7:A25")
$blankA.Font.Color = 0
$blankA.WrapText = $true
$blankA.HorizontalAlignment = -4131

$blankB = $ws.Range("B7:B25")
$blankB.WrapText = $true
$blankB.HorizontalAlignment = 1

$blankC = $ws.Range("C7:C25")
$blankC.Font.Color = 0
$blankC.WrapText = $true
$blankC.HorizontalAlignment = -4131
$blankC.NumberFormat = "#,##0"

# ------------------------------------------------------------------
# Sheet dimensions / frozen pane stay as-is (A1:C25, freeze at row 3)
# ------------------------------------------------------------------
